{"js": "// 1) Insert a literal tab character between \"[ x ]\" and \"Modifications\" in the\n//    Pilot Period modifications line, keeping everything else in the same run.\n{\n  const body = context.document.body;\n  const results = body.search(\n    \"[ x ]Modifications to the Agreement that apply only to the Pilot Period: {pilot_modifications}\",\n    { matchCase: true }\n  );\n  results.load(\"items,text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const target = results.items[0];\n    target.insertText(\n      \"[ x ]\\tModifications to the Agreement that apply only to the Pilot Period: {pilot_modifications}\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// 2) Remove the 4 empty trailing paragraphs that follow the\n//    \"{payment_display}\" paragraph (stability fix for the fill formatting).\n{\n  const body = context.document.body;\n  const results = body.search(\"{payment_display}\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const anchorParagraph = results.items[0].paragraphs.getFirst();\n    let current = anchorParagraph.getNextOrNullObject();\n    current.load(\"isNullObject,text\");\n    await context.sync();\n\n    let removed = 0;\n    while (!current.isNullObject && current.text === \"\" && removed < 4) {\n      const toDelete = current;\n      current = current.getNextOrNullObject();\n      toDelete.delete();\n      await context.sync();\n      current.load(\"isNullObject,text\");\n      await context.sync();\n      removed++;\n    }\n  }\n}\n", "ps1": "# 1) Insert a literal tab character between \"[ x ]\" and \"Modifications\" in the\n#    Pilot Period modifications line, keeping everything else in the same run.\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$oldText = \"[ x ]Modifications to the Agreement that apply only to the Pilot Period: {pilot_modifications}\"\n$found = $findRange.Find.Execute($oldText)\nif ($found) {\n    $findRange.Text = \"[ x ]`tModifications to the Agreement that apply only to the Pilot Period: {pilot_modifications}\"\n}\n\n# 2) Remove the 4 empty trailing paragraphs that follow the \"{payment_display}\"\n#    paragraph in the Payment Process row (stability fix for the fill formatting).\n$table = $d.Tables(1)\n$targetRow = -1\nFor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $cell = $table.Cell($r, 2)\n    if ($cell.Range.Text -like \"*{payment_display}*\") {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -ne -1) {\n    $cell0 = $table.Cell($targetRow, 2)\n    $pdIndex = -1\n    For ($i = 1; $i -le $cell0.Range.Paragraphs.Count; $i++) {\n        $p = $cell0.Range.Paragraphs($i)\n        if ($p.Range.Text -like \"*{payment_display}*\") {\n            $pdIndex = $i\n            break\n        }\n    }\n\n    if ($pdIndex -ne -1) {\n        $deleteIndex = $pdIndex + 1\n        For ($k = 1; $k -le 4; $k++) {\n            # Re-fetch the cell/paragraph fresh every iteration: once a\n            # paragraph is deleted, the next empty paragraph slides into the\n            # same index, and stale object references misbehave here.\n            $cell = $table.Cell($targetRow, 2)\n            $p = $cell.Range.Paragraphs($deleteIndex)\n            $p.Range.Delete() | Out-Null\n        }\n    }\n}\n"}
